$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.811.53'
$ws.Range("E2").Value = '  +1.41%  '
$ws.Range("D3").Value = '2.044.19'
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("E4").Value = '  +0.09%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '229.42'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("E6").Value = '  +1.19%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '58.27'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +5.27%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +1.56%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0808'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.51%  '
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").Value = '2.346.06'
$ws.Range("E12").Value = '  +0.71%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '14.52'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.60%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '20.85'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.12%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.29'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +2.13%  '
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '2.037.46'
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").Value = '37.743.69'
$ws.Range("E18").Value = '  +1.30%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.18'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.47%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '69.70'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("D21").Value = '0.0₃0833'
$ws.Range("E21").Value = '  +1.20%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '224.26'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("E23").Value = '  -0.01%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.44'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.12%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.25'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.38%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '166.59'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("E27").Value = '  -0.69%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.133'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +3.19%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '19.04'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.02%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.34'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.88%  '
$ws.Range("E31").Value = '  +1.26%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.52'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.78%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '2.09'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +13.39%  '
$ws.Range("E34").Value = '  +2.09%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0612'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.87%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.33'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.87%  '
$ws.Range("E37").Value = '  +8.62%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.29'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +3.98%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0218'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.481.44'
$ws.Range("E41").Value = '  +0.10%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '96.75'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.24%  '
$ws.Range("E43").Value = '  +2.20%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.0932'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.38%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '16.60'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.13'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.59%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '4.10'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +14.37%  '
$ws.Range("E48").Value = '  -0.60%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.95'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.74%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '6.98'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -4.00%  '
$ws.Range("D51").Value = '2.233.17'
$ws.Range("E51").Value = '  +0.86%  '
